# Updates cryptos list data (Price / Volume(1h) columns, and row 35/36 coin swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Formula = "'68.953.46"
$ws.Range("E2").Formula = "'  +2.46%  "

$ws.Range("D3").Formula = "'3.731.67"
$ws.Range("E3").Formula = "'  +0.59%  "

$ws.Range("E4").Formula = "'  +0.03%  "

$ws.Range("D5").Formula = "'601.40"
$ws.Range("E5").Formula = "'  +0.67%  "

$ws.Range("D6").Formula = "'167.28"
$ws.Range("E6").Formula = "'  +0.14%  "

$ws.Range("D7").Formula = "'3.727.37"
$ws.Range("E7").Formula = "'  +0.56%  "

$ws.Range("E8").Formula = "'  -0.07%  "

$ws.Range("D9").Formula = "'0.535"
$ws.Range("E9").Formula = "'  -0.02%  "

$ws.Range("D10").Formula = "'0.165"
$ws.Range("E10").Formula = "'  +0.43%  "

$ws.Range("D11").Formula = "'6.41"
$ws.Range("E11").Formula = "'  +3.50%  "

$ws.Range("D12").Formula = "'0.459"
$ws.Range("E12").Formula = "'  -0.44%  "

$ws.Range("D13").Formula = "'37.83"
$ws.Range("E13").Formula = "'  +0.12%  "

$ws.Range("D14").Formula = "'0.0000245"
$ws.Range("E14").Formula = "'  +1.22%  "

$ws.Range("D15").Formula = "'4.364.89"
$ws.Range("E15").Formula = "'  +0.84%  "

$ws.Range("D16").Formula = "'3.741.89"
$ws.Range("E16").Formula = "'  +0.82%  "

$ws.Range("D17").Formula = "'69.115.03"
$ws.Range("E17").Formula = "'  +2.62%  "

$ws.Range("D18").Formula = "'7.27"
$ws.Range("E18").Formula = "'  +0.13%  "

$ws.Range("E19").Formula = "'  -0.89%  "

$ws.Range("D20").Formula = "'16.95"
$ws.Range("E20").Formula = "'  -2.69%  "

$ws.Range("D21").Formula = "'10.78"
$ws.Range("E21").Formula = "'  +16.28%  "

$ws.Range("D22").Formula = "'492.31"
$ws.Range("E22").Formula = "'  +1.03%  "

$ws.Range("D23").Formula = "'0.722"
$ws.Range("E23").Formula = "'  -0.69%  "

$ws.Range("D24").Formula = "'0.0000149"
$ws.Range("E24").Formula = "'  +4.47%  "

$ws.Range("D25").Formula = "'84.64"
$ws.Range("E25").Formula = "'  -0.42%  "

$ws.Range("D26").Formula = "'2.30"
$ws.Range("E26").Formula = "'  +0.30%  "

$ws.Range("D27").Formula = "'12.22"
$ws.Range("E27").Formula = "'  +0.17%  "

$ws.Range("D28").Formula = "'10.10"
$ws.Range("E28").Formula = "'  +0.76%  "

$ws.Range("E29").Formula = "'  +0.01%  "

$ws.Range("D30").Formula = "'2.98"
$ws.Range("E30").Formula = "'  +2.20%  "

$ws.Range("D31").Formula = "'2.50"
$ws.Range("E31").Formula = "'  +5.85%  "

$ws.Range("D32").Formula = "'8.05"
$ws.Range("E32").Formula = "'  +4.78%  "

$ws.Range("D33").Formula = "'31.46"
$ws.Range("E33").Formula = "'  +0.21%  "

$ws.Range("D34").Formula = "'3.887.01"
$ws.Range("E34").Formula = "'  +0.92%  "

$ws.Range("B35").Value = "RenzoRestakedETH"
$ws.Range("C35").Value = "https://coinranking.com/coin/lKlJ_MC5M+renzorestakedeth-ezeth"
$ws.Range("D35").Formula = "'3.674.42"
$ws.Range("E35").Formula = "'  +0.62%  "

$ws.Range("B36").Value = "Hedera"
$ws.Range("C36").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D36").Formula = "'0.108"
$ws.Range("E36").Formula = "'  -0.19%  "

$ws.Range("D37").Formula = "'1.00"
$ws.Range("E37").Formula = "'  +0.14%  "

$ws.Range("D38").Formula = "'1.01"
$ws.Range("E38").Formula = "'  +1.38%  "

$ws.Range("D39").Formula = "'5.85"
$ws.Range("E39").Formula = "'  +0.38%  "

$ws.Range("E40").Formula = "'  +1.15%  "

$ws.Range("D41").Formula = "'0.322"
$ws.Range("E41").Formula = "'  +0.23%  "

$ws.Range("D42").Formula = "'2.97"
$ws.Range("E42").Formula = "'  +5.07%  "

$ws.Range("D43").Formula = "'430.89"
$ws.Range("E43").Formula = "'  +0.90%  "

$ws.Range("D44").Formula = "'48.56"
$ws.Range("E44").Formula = "'  -0.28%  "

$ws.Range("D45").Formula = "'1.97"
$ws.Range("E45").Formula = "'  +2.18%  "

$ws.Range("D46").Formula = "'8.44"
$ws.Range("E46").Formula = "'  -0.11%  "

$ws.Range("E47").Formula = "'  -0.01%  "

$ws.Range("D48").Formula = "'40.10"
$ws.Range("E48").Formula = "'  -0.84%  "

$ws.Range("D49").Formula = "'140.33"
$ws.Range("E49").Formula = "'  -0.05%  "

$ws.Range("D50").Formula = "'2.765.46"
$ws.Range("E50").Formula = "'  +0.41%  "

$ws.Range("D51").Formula = "'0.0352"
$ws.Range("E51").Formula = "'  +0.51%  "
